$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 9,28
$data[0,0] = "M1"
$data[0,1] = "A1"
$data[0,2] = "DO"
$data[0,3] = "M1"
$data[0,4] = "PH"
$data[0,5] = "M3"
$data[0,6] = "PH"
$data[0,7] = "A1"
$data[0,8] = "DO"
$data[0,9] = "PH"
$data[0,10] = "A1"
$data[0,11] = "M3"
$data[0,12] = "M1"
$data[0,13] = "PH"
$data[0,14] = "A1"
$data[0,15] = "A1"
$data[0,16] = "DO"
$data[0,17] = "M3"
$data[0,18] = "PH"
$data[0,19] = "M1"
$data[0,20] = "PH"
$data[0,21] = "M3"
$data[0,22] = "M1"
$data[0,23] = "DO"
$data[0,24] = "A1"
$data[0,25] = "PH"
$data[0,26] = "M1"
$data[0,27] = "M1"
$data[1,0] = "A1"
$data[1,1] = "A1"
$data[1,2] = "A1"
$data[1,3] = "DO"
$data[1,4] = "PH"
$data[1,5] = "M3"
$data[1,6] = "PH"
$data[1,7] = "M1"
$data[1,8] = "A1"
$data[1,9] = "A1"
$data[1,10] = "PH"
$data[1,11] = "M3"
$data[1,12] = "M2"
$data[1,13] = "DO"
$data[1,14] = "A1"
$data[1,15] = "PH"
$data[1,16] = "PH"
$data[1,17] = "M3"
$data[1,18] = "DO"
$data[1,19] = "M1"
$data[1,20] = "M2"
$data[1,21] = "M1"
$data[1,22] = "DO"
$data[1,23] = "M2"
$data[1,24] = "M3"
$data[1,25] = "M2"
$data[1,26] = "M2"
$data[1,27] = "PH"
$data[2,0] = "DO"
$data[2,1] = "M1"
$data[2,2] = "M1"
$data[2,3] = "M1"
$data[2,4] = "M3"
$data[2,5] = "M1"
$data[2,6] = "M1"
$data[2,7] = "M1"
$data[2,8] = "M1"
$data[2,9] = "DO"
$data[2,10] = "M1"
$data[2,11] = "M3"
$data[2,12] = "A1"
$data[2,13] = "M3"
$data[2,14] = "M3"
$data[2,15] = "M1"
$data[2,16] = "A1"
$data[2,17] = "A1"
$data[2,18] = "M1"
$data[2,19] = "M1"
$data[2,20] = "DO"
$data[2,21] = "M1"
$data[2,22] = "PH"
$data[2,23] = "A1"
$data[2,24] = "DO"
$data[2,25] = "PH"
$data[2,26] = "M3"
$data[2,27] = "M1"
$data[3,0] = "DO"
$data[3,1] = "M3"
$data[3,2] = "M1"
$data[3,3] = "M1"
$data[3,4] = "M1"
$data[3,5] = "M1"
$data[3,6] = "PH"
$data[3,7] = "DO"
$data[3,8] = "M2"
$data[3,9] = "M1"
$data[3,10] = "M3"
$data[3,11] = "PH"
$data[3,12] = "M2"
$data[3,13] = "PH"
$data[3,14] = "M1"
$data[3,15] = "DO"
$data[3,16] = "M1"
$data[3,17] = "M1"
$data[3,18] = "M1"
$data[3,19] = "M3"
$data[3,20] = "PH"
$data[3,21] = "A1"
$data[3,22] = "M1"
$data[3,23] = "M2"
$data[3,24] = "M2"
$data[3,25] = "M3"
$data[3,26] = "M3"
$data[3,27] = "DO"
$data[4,0] = "M1"
$data[4,1] = "DO"
$data[4,2] = "A2"
$data[4,3] = "A1"
$data[4,4] = "A2"
$data[4,5] = "M3"
$data[4,6] = "A1"
$data[4,7] = "A1"
$data[4,8] = "A1"
$data[4,9] = "M3"
$data[4,10] = "DO"
$data[4,11] = "A1"
$data[4,12] = "M1"
$data[4,13] = "A2"
$data[4,14] = "A1"
$data[4,15] = "A2"
$data[4,16] = "M3"
$data[4,17] = "DO"
$data[4,18] = "M1"
$data[4,19] = "M2"
$data[4,20] = "M1"
$data[4,21] = "A1"
$data[4,22] = "A1"
$data[4,23] = "DO"
$data[4,24] = "M1"
$data[4,25] = "A1"
$data[4,26] = "M3"
$data[4,27] = "A1"
$data[5,0] = "A1"
$data[5,1] = "M3"
$data[5,2] = "DO"
$data[5,3] = "A1"
$data[5,4] = "A1"
$data[5,5] = "A1"
$data[5,6] = "A1"
$data[5,7] = "DO"
$data[5,8] = "M3"
$data[5,9] = "A1"
$data[5,10] = "A1"
$data[5,11] = "A1"
$data[5,12] = "A1"
$data[5,13] = "A1"
$data[5,14] = "DO"
$data[5,15] = "M3"
$data[5,16] = "A1"
$data[5,17] = "A1"
$data[5,18] = "A1"
$data[5,19] = "A1"
$data[5,20] = "A1"
$data[5,21] = "DO"
$data[5,22] = "M3"
$data[5,23] = "A1"
$data[5,24] = "A1"
$data[5,25] = "A1"
$data[5,26] = "A1"
$data[5,27] = "A1"
$data[6,0] = "DO"
$data[6,1] = "M1"
$data[6,2] = "A2"
$data[6,3] = "A1"
$data[6,4] = "A2"
$data[6,5] = "M3"
$data[6,6] = "M1"
$data[6,7] = "A1"
$data[6,8] = "DO"
$data[6,9] = "A2"
$data[6,10] = "A2"
$data[6,11] = "A2"
$data[6,12] = "M2"
$data[6,13] = "M3"
$data[6,14] = "DO"
$data[6,15] = "A1"
$data[6,16] = "A1"
$data[6,17] = "A1"
$data[6,18] = "A1"
$data[6,19] = "M2"
$data[6,20] = "M3"
$data[6,21] = "M3"
$data[6,22] = "A2"
$data[6,23] = "A1"
$data[6,24] = "A1"
$data[6,25] = "A1"
$data[6,26] = "M2"
$data[6,27] = "DO"
$data[7,0] = "M1"
$data[7,1] = "DO"
$data[7,2] = "M1"
$data[7,3] = "M3"
$data[7,4] = "M2"
$data[7,5] = "A2"
$data[7,6] = "A1"
$data[7,7] = "DO"
$data[7,8] = "M1"
$data[7,9] = "M2"
$data[7,10] = "M2"
$data[7,11] = "M3"
$data[7,12] = "A2"
$data[7,13] = "A2"
$data[7,14] = "DO"
$data[7,15] = "M1"
$data[7,16] = "M1"
$data[7,17] = "M1"
$data[7,18] = "M3"
$data[7,19] = "A2"
$data[7,20] = "A1"
$data[7,21] = "DO"
$data[7,22] = "M1"
$data[7,23] = "M3"
$data[7,24] = "M2"
$data[7,25] = "M1"
$data[7,26] = "A2"
$data[7,27] = "A2"
$data[8,0] = "A1"
$data[8,1] = "A1"
$data[8,2] = "M1"
$data[8,3] = "DO"
$data[8,4] = "M3"
$data[8,5] = "A2"
$data[8,6] = "M2"
$data[8,7] = "M1"
$data[8,8] = "A2"
$data[8,9] = "M1"
$data[8,10] = "M2"
$data[8,11] = "DO"
$data[8,12] = "M3"
$data[8,13] = "M1"
$data[8,14] = "M3"
$data[8,15] = "M1"
$data[8,16] = "M1"
$data[8,17] = "DO"
$data[8,18] = "A2"
$data[8,19] = "A2"
$data[8,20] = "A2"
$data[8,21] = "A1"
$data[8,22] = "A2"
$data[8,23] = "M1"
$data[8,24] = "DO"
$data[8,25] = "M3"
$data[8,26] = "A2"
$data[8,27] = "M1"
$ws.Range("B2:AC10").Value = $data
